$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 22
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 724
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 61
$ws.Range("L2").Value = "baby"
$ws.Range("M2").Value = "masculin"
$ws.Range("N2").Value = "code"
$ws.Range("O2").Value = "test"
